## Apply "taking latest changes and adding new profile implemented scripts"
## - Switches Runmode (col D) from Y to N for all currently existing test rows (2-25)
## - Switches row 10's Result (col E) from PASS to SKIP
## - Appends two new test-case rows (26, 27) for the new Profile Primary Institution
##   type-ahead scripts

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Rows 26/27 reuse the same formatting as row 25 (thin border all round, plain font,
# highlighted "JIRA ID" column B) -- copy it down before filling in the new values.
$ws.Range("A25:E25").Copy()
$ws.Range("A26:E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New row 26: ProfilePrimaryInstitutionTypeAheadTest (written first so the shared
#     string table gains the new strings -- including "N" -- in source order) ---
$ws.Cells.Item(26, 1).Value = "ProfilePrimaryInstitutionTypeAheadTest"
$ws.Cells.Item(26, 2).Value = "TBD"
$ws.Cells.Item(26, 3).Value = "Verify that user is able to add 'primary institution' using type ahead"
$ws.Cells.Item(26, 4).Value = "N"
$ws.Cells.Item(26, 5).Value = "SKIP"

# --- Update Runmode (column D) for existing rows 2..25: Y -> N ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 4).Value = "N"
}

# Row 10's Result (column E) changes from PASS to SKIP
$ws.Cells.Item(10, 5).Value = "SKIP"

# --- New row 27: ProfilePrimaryInstitutionTypeAheadMinCharTest (mixed formatting description) ---
$ws.Cells.Item(27, 1).Value = "ProfilePrimaryInstitutionTypeAheadMinCharTest"
$ws.Cells.Item(27, 2).Value = "TBD"

$descCell = $ws.Cells.Item(27, 3)
$fullText = "Verify that  'primary institution' type ahead options should display while enter min 2 characters"
$descCell.Value = $fullText

# Italicize the "primary institution' " portion (characters 15-35, 1-based)
$italicRun = $descCell.Characters(15, 21)
$italicRun.Font.Italic = $true
$italicRun.Font.Size = 11
$italicRun.Font.Name = "Calibri"

# Remaining run keeps the normal (non-italic) font explicitly
$restRun = $descCell.Characters(36, 62)
$restRun.Font.Size = 11
$restRun.Font.Name = "Calibri"

$ws.Cells.Item(27, 4).Value = "Y"
$ws.Cells.Item(27, 5).Value = "PASS"

# Selection / active cell moves to C19 after editing
$ws.Range("C19").Select()
